# Apply "update RNAseq report and doc to Aladdin branding" edit.
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The second (duplicate) "preseq" hyperlink was left empty
#    (<w:hyperlink r:id="rId14" .../> with no run/text inside it).
#    Populate it with the same display text as its neighbour so the
#    hyperlink actually shows/works. This is done FIRST, before any
#    other content edit, because a still-empty hyperlink element gets
#    pruned away as soon as its host paragraph is touched by an
#    unrelated edit.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -eq "https://github.com/smithlabcode/preseq" -and $h.TextToDisplay -eq "") {
        $h.TextToDisplay = "https://github.com/smithlabcode/preseq"
    }
}

# ------------------------------------------------------------------
# 2) Rebrand the three "Zymo Research" mentions to "Aladdin" (the
#    company name changed; the hyphenated "Zymo-Research" inside the
#    old github URL text is intentionally left untouched, since the
#    search string below requires a literal space after "Zymo").
# ------------------------------------------------------------------
$d.Content.Find.Execute("Zymo Research ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Aladdin ", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Update the citation line: "Aladdin Open Bioinformatics, 2021."
#    -> "Aladdin Bioinformatics Platform, 2022."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Aladdin Open Bioinformatics, 2021.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Aladdin Bioinformatics Platform, 2022.", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Remove the stray "_GoBack" bookmark that trailed the citation.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
